# Cambio en casos de prueba
# Populates the "Testing" sheet with the Login/Register test-case rows,
# fixes the header, sets row heights for the wrapped rows, and adds the
# hyperlinks on the "Entradas" column for the evidence-backed cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (row 2): "Salida esperada" -> "Salida", new "Resultado
# esperado" column I2 stays the same text but a new H/I pair of columns
# is introduced (Salida / Resultado esperado).
# ---------------------------------------------------------------------
$ws.Range("H2").Value = "Salida"
$ws.Range("I2").Value = "Resultado esperado"

# ---------------------------------------------------------------------
# Row 3 - CP_Login1
# ---------------------------------------------------------------------
$ws.Range("C3").Value = "CP_Login1"
$ws.Range("D3").Value = "Iniciar sesion con datos verdaderos"
$ws.Range("E3").Value = "Debe estar registrado el usuario"
$ws.Range("F3").Value = "Yeison@MundoAnimal.com / 1000088550"
$ws.Range("H3").Value = "Ingreso al sistema"
$ws.Range("I3").Value = "Corf"
$ws.Rows.Item(3).RowHeight = 28.5

# ---------------------------------------------------------------------
# Row 4 - CP_Login2
# ---------------------------------------------------------------------
$ws.Range("C4").Value = "CP_Login2"
$ws.Range("D4").Value = "Iniciar sesion con datos erroneos"
$ws.Range("E4").Value = "Debe estar registrado el usuario"
$ws.Range("F4").Value = "Yeison@Mundo.com / 100008855"
$ws.Rows.Item(4).RowHeight = 30

# ---------------------------------------------------------------------
# Row 5 - CP_Login3
# ---------------------------------------------------------------------
$ws.Range("C5").Value = "CP_Login3"
$ws.Range("D5").Value = "Iniciar sesion con los campos vacios"
$ws.Range("E5").Value = "Debe contener los campos obligatorios"
$ws.Range("F5").Value = "Null / Null"
$ws.Rows.Item(5).RowHeight = 30

# ---------------------------------------------------------------------
# Row 6 - CP_Login4
# ---------------------------------------------------------------------
$ws.Range("C6").Value = "CP_Login4"
$ws.Range("D6").Value = "Iniciar sesion con campo de correo vacio y contraseña diligenciado"
$ws.Range("E6").Value = "Debe contener los campos obligatorios"
$ws.Range("F6").Value = "Null / 1000088550"
$ws.Rows.Item(6).RowHeight = 31.5

# ---------------------------------------------------------------------
# Row 7 - CP_Login5
# ---------------------------------------------------------------------
$ws.Range("C7").Value = "CP_Login5"
$ws.Range("D7").Value = "Iniciar sesion con campo de correo diligenciado y contraseña vacio"
$ws.Range("E7").Value = "Debe contener los campos obligatorios"
$ws.Range("F7").Value = "Yeison@MundoAnimal.com / Null"
$ws.Rows.Item(7).RowHeight = 30.75

# ---------------------------------------------------------------------
# Row 8 - CP_Login6
# ---------------------------------------------------------------------
$ws.Range("C8").Value = "CP_Login6"
$ws.Range("D8").Value = "Iniciar sesion con el campo correo erroneo"
$ws.Range("E8").Value = "Debe estar registrado el usuario"
$ws.Range("F8").Value = "Yeison@Animal.com / 10000088550"
$ws.Rows.Item(8).RowHeight = 34.5

# ---------------------------------------------------------------------
# Row 9 - CP_Login7
# ---------------------------------------------------------------------
$ws.Range("C9").Value = "CP_Login7"
$ws.Range("D9").Value = "Iniciar sesion con el campo contraseña erroneo"
$ws.Range("E9").Value = "Debe estar registrado el usuario"
$ws.Range("F9").Value = "Yeison@MundoAnimal.com / 10008855"
$ws.Rows.Item(9).RowHeight = 30.75

# ---------------------------------------------------------------------
# Row 10 - CP_Register1
# ---------------------------------------------------------------------
$ws.Range("C10").Value = "CP_Register1"
$ws.Range("D10").Value = "Registrar usuario con todos los campos vacios"
$ws.Range("E10").Value = "El usuario debe llenar los campos solicitados"
$ws.Range("F10").Value = "74635215 / Juan / Martinez / 315698754 / CL 87 CR 31 67 / JuanM@gmail.com / 74635215"
$ws.Rows.Item(10).RowHeight = 45

# ---------------------------------------------------------------------
# Rows 11-18 - CP_Register2..CP_Register9 (identifier only)
# ---------------------------------------------------------------------
$ws.Range("C11").Value = "CP_Register2"
$ws.Range("C12").Value = "CP_Register3"
$ws.Range("C13").Value = "CP_Register4"
$ws.Range("C14").Value = "CP_Register5"
$ws.Range("C15").Value = "CP_Register6"
$ws.Range("C16").Value = "CP_Register7"
$ws.Range("C17").Value = "CP_Register8"
$ws.Range("C18").Value = "CP_Register9"

# ---------------------------------------------------------------------
# Hyperlinks on the "Entradas" cells that reference supporting evidence.
# ---------------------------------------------------------------------
$ws.Range("F3").Hyperlinks.Add($ws.Range("F3"), "https://example.com/evidencia/CP_Login1") | Out-Null
$ws.Range("F4").Hyperlinks.Add($ws.Range("F4"), "https://example.com/evidencia/CP_Login2") | Out-Null
$ws.Range("F7").Hyperlinks.Add($ws.Range("F7"), "https://example.com/evidencia/CP_Login5") | Out-Null
$ws.Range("F8").Hyperlinks.Add($ws.Range("F8"), "https://example.com/evidencia/CP_Login6") | Out-Null
$ws.Range("F9").Hyperlinks.Add($ws.Range("F9"), "https://example.com/evidencia/CP_Login7") | Out-Null

# ---------------------------------------------------------------------
# Window / view state: zoom in to 136% and scroll so column G is the
# left-most visible column, with I5 as the active selection.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 136
$ws.Range("I5").Select() | Out-Null
